$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells to remain text (avoid Excel auto-numeric conversion),
# matching the original inlineStr cell type in the workbook.
$priceCells = @("D2","D3","D5","D6","D7","D9","D11","D13","D14","D15","D17","D18","D19","D20","D22","D23","D24","D25","D26","D29","D30","D31","D32","D34","D36","D38","D39","D40","D41","D42","D44","D46","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "63.175.61"
$ws.Range("E2").Value = "  -0.73%  "

# Row 3
$ws.Range("D3").Value = "3.236.51"
$ws.Range("E3").Value = "  -0.74%  "

# Row 4
$ws.Range("E4").Value = "  -0.43%  "

# Row 5
$ws.Range("D5").Value = "529.95"
$ws.Range("E5").Value = "  +3.88%  "

# Row 6
$ws.Range("D6").Value = "172.32"
$ws.Range("E6").Value = "  -1.69%  "

# Row 7
$ws.Range("D7").Value = "0.595"
$ws.Range("E7").Value = "  +1.73%  "

# Row 8
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("D9").Value = "3.235.09"
$ws.Range("E9").Value = "  -0.50%  "

# Row 10
$ws.Range("E10").Value = "  -0.08%  "

# Row 11
$ws.Range("D11").Value = "53.23"
$ws.Range("E11").Value = "  -5.82%  "

# Row 12
$ws.Range("E12").Value = "  +4.04%  "

# Row 13
$ws.Range("D13").Value = "0.0000254"
$ws.Range("E13").Value = "  +1.32%  "

# Row 14
$ws.Range("D14").Value = "9.13"
$ws.Range("E14").Value = "  +2.40%  "

# Row 15
$ws.Range("D15").Value = "3.752.30"
$ws.Range("E15").Value = "  -1.34%  "

# Row 16
$ws.Range("E16").Value = "  -1.67%  "

# Row 17
$ws.Range("D17").Value = "3.233.98"
$ws.Range("E17").Value = "  -1.11%  "

# Row 18
$ws.Range("D18").Value = "63.065.85"
$ws.Range("E18").Value = "  -0.63%  "

# Row 19
$ws.Range("D19").Value = "17.21"
$ws.Range("E19").Value = "  +2.13%  "

# Row 20
$ws.Range("D20").Value = "11.08"
$ws.Range("E20").Value = "  +4.15%  "

# Row 21
$ws.Range("E21").Value = "  +3.82%  "

# Row 22
$ws.Range("D22").Value = "366.48"
$ws.Range("E22").Value = "  +0.71%  "

# Row 23
$ws.Range("D23").Value = "3.78"
$ws.Range("E23").Value = "  +5.16%  "

# Row 24
$ws.Range("B24").Value = "RenderToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D24").Value = "11.27"
$ws.Range("E24").Value = "  +5.85%  "

# Row 25
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "81.27"
$ws.Range("E25").Value = "  +2.90%  "

# Row 26
$ws.Range("D26").Value = "4.01"
$ws.Range("E26").Value = "  +6.74%  "

# Row 27
$ws.Range("E27").Value = "  +2.68%  "

# Row 28
$ws.Range("E28").Value = "  +2.17%  "

# Row 29
$ws.Range("D29").Value = "11.30"
$ws.Range("E29").Value = "  +1.69%  "

# Row 30
$ws.Range("D30").Value = "8.18"
$ws.Range("E30").Value = "  +0.12%  "

# Row 31
$ws.Range("D31").Value = "28.51"
$ws.Range("E31").Value = "  +1.70%  "

# Row 32
$ws.Range("D32").Value = "638.10"
$ws.Range("E32").Value = "  +1.34%  "

# Row 33
$ws.Range("E33").Value = "  -1.44%  "

# Row 34
$ws.Range("D34").Value = "11.24"
$ws.Range("E34").Value = "  +3.30%  "

# Row 35
$ws.Range("E35").Value = "  +4.72%  "

# Row 36
$ws.Range("D36").Value = "56.93"
$ws.Range("E36").Value = "  -3.46%  "

# Row 37
$ws.Range("E37").Value = "  -0.01%  "

# Row 38
$ws.Range("D38").Value = "36.82"
$ws.Range("E38").Value = "  +4.67%  "

# Row 39
$ws.Range("D39").Value = "0.378"
$ws.Range("E39").Value = "  +2.55%  "

# Row 40
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.13%  "

# Row 41
$ws.Range("D41").Value = "0.0₃0711"
$ws.Range("E41").Value = "  +13.16%  "

# Row 42
$ws.Range("D42").Value = "2.62"
$ws.Range("E42").Value = "  +15.69%  "

# Row 43
$ws.Range("E43").Value = "  +2.18%  "

# Row 44
$ws.Range("D44").Value = "2.886.94"
$ws.Range("E44").Value = "  +1.84%  "

# Row 45
$ws.Range("E45").Value = "  +12.17%  "

# Row 46
$ws.Range("D46").Value = "2.69"
$ws.Range("E46").Value = "  +4.72%  "

# Row 47
$ws.Range("E47").Value = "  +4.44%  "

# Row 48
$ws.Range("E48").Value = "  -1.40%  "

# Row 49
$ws.Range("E49").Value = "  +7.34%  "

# Row 50
$ws.Range("E50").Value = "  +2.14%  "

# Row 51
$ws.Range("D51").Value = "134.67"
$ws.Range("E51").Value = "  +2.08%  "

# Clear the temporary text-number-format styling from Price cells so no new
# style index is introduced (the source cells carried no explicit style).
foreach ($addr in $priceCells) {
    $ws.Range($addr).ClearFormats()
}
